# Auto commit update: refresh Metrics figures and dependent "today" sheet
# formulas, then restore each sheet's last-used selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Metrics" sheet - update the raw figures in column B (rows 2-13).
#    Every downstream formula (on "today") references these cells, so they
#    will recalculate automatically once these values change.
# ---------------------------------------------------------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 94403.78
$metrics.Range("B3").Value  = 81254.180000000008
$metrics.Range("B4").Value  = 28980.479999999996
$metrics.Range("B5").Value  = 3834
$metrics.Range("B6").Value  = 5297110.8900000006
$metrics.Range("B7").Value  = 4481607.1400000006
$metrics.Range("B8").Value  = 1560937.3600000003
$metrics.Range("B9").Value  = 206541
$metrics.Range("B10").Value = 33762491.879999995
$metrics.Range("B11").Value = 31756882.300000004
$metrics.Range("B12").Value = 11842659.399999997
$metrics.Range("B13").Value = 1304171

# Restore this sheet's recorded selection (was E28, now E19).
$metrics.Activate()
$metrics.Range("E19").Select()

# ---------------------------------------------------------------------------
# 2. "today" sheet - its B/E/F columns are formulas that reference Metrics,
#    so they recalc on their own. Only the selection needs to move (was
#    E10, now F7). The A1 TODAY()-1 cell is volatile and recalculates with
#    the host clock, so it needs no manual update here.
# ---------------------------------------------------------------------------
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("F7").Select()
